# Make the line spacing default in the Normal style match the input
# file's: change it from 1.5 lines to double spacing (the font stays
# untouched, per the commit message).
$d = $word.ActiveDocument

$normal = $d.Styles("Normal")
$normal.ParagraphFormat.LineSpacingRule = 2   # wdLineSpaceDouble -> w:line="480" w:lineRule="auto"
